# "Mark latrine questions as required"
#
# The "survey" sheet holds the ODK XLSForm-style question definitions, one
# row per clause/question (columns: clause, condition, type, values_list,
# name, display.prompt, constraint, inputAttributes.min,
# inputAttributes.step). This adds a new "required" column (J) and flags
# the relevant questions as required:
#   - row 4: latrine_type            -> required = 1
#   - row 6: latrine_shared          -> required = 1
#   - row 8: latrine_shared_num_hh   -> required = the same condition that
#            already guards whether the question is shown at all
#            (selected(data('latrine_shared'), 'yes'))
#
# The "settings" sheet's form_version is bumped to reflect the new
# form revision, and sheet focus moves from "survey" to "settings".

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$survey.Range("J1").Value = "required"
$survey.Range("J4").Value = 1
$survey.Range("J6").Value = 1
$survey.Range("J8").Value = "selected(data('latrine_shared'), 'yes')"
$survey.Range("J9").Select() | Out-Null

$settings = $wb.Worksheets.Item("settings")
$settings.Range("B3").Value = 20210421001
$settings.Activate() | Out-Null
